$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1-3) Collapse the three "gramStart/gramEnd"-split paragraphs back into a
#      single run each (the proofing-error wrapper runs go away, the text is
#      unchanged once concatenated).
# ---------------------------------------------------------------------------
$p1 = "Obrigatório o uso do cinturão de segurança em atividades com altura igual ou maior que dois metros, ou quando o risco da atividade assim determinar;"
$d.Content.Find.Execute($p1, $true, $false, $false, $false, $false, $true, 1, $false, $p1, 2) | Out-Null

$p2 = "Só execute serviços ou opere máquinas se estiver devidamente habilitado e autorizado, quando não souber ou tiver dúvidas sobre algum serviço, pergunte ao seu superior antes do início ou durante a realização do mesmo, para prevenir-se contra possíveis acidentes;"
$d.Content.Find.Execute($p2, $true, $false, $false, $false, $false, $true, 1, $false, $p2, 2) | Out-Null

$p3 = "Em caso de eventuais Acidentes do Trabalho, o funcionário deve de imediato comunicar a Segurança do Trabalho, chefia e/ou responsável, para que o mesmo receba os primeiros socorros e seja feita a abertura da Comunicação de Acidentes de Trabalho – CAT; encaminhando à vítima ao ambulatório e/ou o hospital/posto de atendimento mais próximo de acordo com o Plano de Emergência do parque eólico."
$d.Content.Find.Execute($p3, $true, $false, $false, $false, $false, $true, 1, $false, $p3, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Signature block: swap the signatory name/title/registration.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("LEONARDO SILVERIO FERREIRA", $true, $false, $false, $false, $false, $true, 1, $false, "BRUNA PETRONI CEZARIO", 2) | Out-Null

$d.Content.Find.Execute("Técnico Segurança do Trabalho", $true, $false, $false, $false, $false, $true, 1, $false, "Engenheira de Segurança do Trabalho", 2) | Out-Null

$d.Content.Find.Execute("MTE/RN: 1360", $true, $false, $false, $false, $false, $true, 1, $false, "CREA-RN: 2122993685", 2) | Out-Null

Write-Output "done"
